$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- column width for column B (best-effort; engine quantizes to 1/6 char) ---
$ws.Range("B:B").ColumnWidth = 11.6

# --- header row 37 ---
$ws.Range("A37").Value = "N"
$ws.Range("B37").Value = "u-uh"
$ws.Range("C37").Value = "time"
$ws.Range("D37").Value = "iter"

# --- row 38 ---
$ws.Range("A38").Value = 32
$ws.Range("B38").Formula = "=0.00301277655464715"
$ws.Range("C38").Value = 0.0122599601745605
$ws.Range("C38").NumberFormat = "0.00E+00"
$ws.Range("D38").Value = 48
$ws.Range("E38").Value = "Ratio"

# --- row 39 ---
$ws.Range("A39").Value = 64
$ws.Range("B39").Value = 0.000778106197006023
$ws.Range("B39").NumberFormat = "0.00E+00"
$ws.Range("C39").Value = 0.00298619270324707
$ws.Range("C39").NumberFormat = "0.00E+00"
$ws.Range("D39").Value = 96
$ws.Range("E39").Value = "N/A"

# --- row 40 ---
$ws.Range("A40").Value = 128
$ws.Range("B40").Value = 0.000197648347696088
$ws.Range("B40").NumberFormat = "0.00E+00"
$ws.Range("C40").Value = 0.032128095626831
$ws.Range("C40").NumberFormat = "0.00E+00"
$ws.Range("D40").Value = 192
$ws.Range("E40").Formula = "=B39/B40"
$ws.Range("E40").NumberFormat = "0.00E+00"

# --- row 41 ---
$ws.Range("A41").Value = 256
$ws.Range("B41").Value = 0.0000497974385832433
$ws.Range("B41").NumberFormat = "0.00E+00"
$ws.Range("C41").Value = 0.168781995773315
$ws.Range("C41").NumberFormat = "0.00E+00"
$ws.Range("D41").Value = 387
$ws.Range("E41").Formula = "=B40/B41"
$ws.Range("E41").NumberFormat = "0.00E+00"

# --- row 42 ---
$ws.Range("A42").Value = 512
$ws.Range("B42").Value = 0.0000124942786493553
$ws.Range("B42").NumberFormat = "0.00E+00"
$ws.Range("C42").Value = 1.24350500106811
$ws.Range("C42").NumberFormat = "0.00E+00"
$ws.Range("D42").Value = 783
$ws.Range("E42").Formula = "=B41/B42"
$ws.Range("E42").NumberFormat = "0.00E+00"

# --- row 43 ---
$ws.Range("A43").Value = 1024
$ws.Range("B43").Value = 0.00000312660644896656
$ws.Range("B43").NumberFormat = "0.00E+00"
$ws.Range("C43").Value = 12.4322278499603
$ws.Range("C43").NumberFormat = "0.00E+00"
$ws.Range("D43").Value = 1581
$ws.Range("E43").Formula = "=B42/B43"
$ws.Range("E43").NumberFormat = "0.00E+00"

# --- row 44 ---
$ws.Range("A44").Value = 2048
$ws.Range("B44").Value = 0.000000780193853411326
$ws.Range("B44").NumberFormat = "0.00E+00"
$ws.Range("C44").Value = 124.338080883026
$ws.Range("C44").NumberFormat = "0.00E+00"
$ws.Range("D44").Value = 3192
$ws.Range("E44").Formula = "=B43/B44"
$ws.Range("E44").NumberFormat = "0.00E+00"

# --- row 45 (B45 empty -> #DIV/0!) ---
$ws.Range("A45").Value = 4096
$ws.Range("E45").Formula = "=B44/B45"
$ws.Range("E45").NumberFormat = "0.00E+00"

# --- row 46 ---
$ws.Range("A46").Value = 8162
$ws.Range("E46").Formula = "=B45/B46"
$ws.Range("E46").NumberFormat = "0.00E+00"

# --- row 47 ---
$ws.Range("A47").Value = 16384
$ws.Range("E47").Formula = "=B46/B47"
$ws.Range("E47").NumberFormat = "0.00E+00"

# --- row 48 ---
$ws.Range("A48").Value = 32768
$ws.Range("E48").Formula = "=B47/B48"
$ws.Range("E48").NumberFormat = "0.00E+00"

# --- sheet view: scroll + selection ---
$ws.Range("C46").Select()

Write-Output "done"
